$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing Excel to keep it as literal text,
# so numeric-looking strings (e.g. "211.46", "2.20") are not silently
# reinterpreted/renormalized as floating point numbers. The temporary
# "@" (Text) number format forces literal storage, then the cell style is
# reset back to Normal so no stray formatting is left behind.
function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Apply updated crypto values scraped on Fri Sep 22 11:59:43 UTC 2023
$ws.Range("D2").Value = "26.679.88"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.599.00"
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue "D5" "211.46"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.17%  "
Set-TextValue "D10" "19.72"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").Value = "1.822.87"
$ws.Range("D13").Value = "1.589.53"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("E15").Value = "  -1.33%  "
Set-TextValue "D16" "64.91"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "26.670.73"
$ws.Range("E18").Value = "  -0.27%  "
Set-TextValue "D19" "210.14"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("E20").Value = "  +0.12%  "
Set-TextValue "D21" "6.79"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  -2.07%  "
$ws.Range("E24").Value = "  +0.54%  "
Set-TextValue "D25" "146.40"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  +0.03%  "
Set-TextValue "D27" "7.20"
$ws.Range("E27").Value = "  -3.74%  "
$ws.Range("E28").Value = "  +2.52%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").Value = "1.296.67"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("E38").Value = "  -1.03%  "
Set-TextValue "D39" "0.844"
$ws.Range("E39").Value = "  +2.72%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D42" "2.20"
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D43" "0.788"
$ws.Range("E43").Value = "  -0.07%  "
Set-TextValue "D44" "63.88"
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("D45").Value = "1.735.64"
$ws.Range("E45").Value = "  +0.22%  "
Set-TextValue "D46" "0.892"
$ws.Range("E46").Value = "  +9.59%  "
Set-TextValue "D47" "90.04"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("E48").Value = "  +0.67%  "
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("E51").Value = "  +1.30%  "
